# Added reading attrition schemes (v0.1.2)
# Adds a new named attrition scheme ("Attrition 1") as row 3 of the
# "Attrition" sheet, and leaves the UI focused on that new row.

$wb = $excel.ActiveWorkbook

$attrition = $wb.Worksheets.Item("Attrition")

# --- Populate the new "Attrition 1" scheme row (row 3) -------------------
$attrition.Range("A3").Value = "Attrition 1"   # Name
$attrition.Range("B3").Value = 12              # Period Length (m)
# C3 already holds the "# nodes" formula; its cached result recalculates
# automatically once the cycle/rate pairs below are populated.
$attrition.Range("D3").Value = 0               # Default
$attrition.Range("E3").Value = 0.02            # cycle 1 rate
$attrition.Range("F3").Value = 5               # cycle 1 period
$attrition.Range("G3").Value = 0.015           # cycle 2 rate
$attrition.Range("H3").Value = 10              # cycle 2 period
$attrition.Range("I3").Value = 0.01            # cycle 3 rate

# --- Make the Attrition sheet the active sheet/selection ------------------
$attrition.Activate() | Out-Null
$attrition.Range("H3").Select() | Out-Null
